# Version 5_0_1 + ShinyFMBN v3 uploaded
# New and larger version of the database + significantly improved version of the app
#
# Updates the "Version_history" sheet:
#  - row 17 (v5): adds the published citation to column C
#  - row 18 (v5.0.1): replaces the short release note / abstract note with the
#    full release note, and updates the abstract note to reference ECCO 2024

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$citation = "Parente, E., Ricciardi, A., 2024. A Comprehensive View of Food Microbiota: Introducing" + [char]10 + "FoodMicrobionet v5. Foods, 13, 1689. https://doi.org/10.3390/foods13111689"
$abstractNote = "Abstract submitted to ECCO 2024"
$releaseNote = "Minimal chances to sample table, added 20 strudies (6 on fungi, 2 on bacteria, 12 on bacteria and fungi) ,removed duplicated samples for ST137"

# Order matters for shared-string slot reuse parity with the source edit:
# citation first, then the abstract note, then the release note.
$ws.Range("C17").Value = $citation
$ws.Range("C18").Value = $abstractNote
$ws.Range("B18").Value = $releaseNote

# Row heights grew to fit the longer wrapped text.
$ws.Rows("17").RowHeight = 75
$ws.Rows("18").RowHeight = 45

# Leave the selection on the cells that were edited last (B18:C18), matching
# the cursor position left behind by the author's edit.
$ws.Range("B18:C18").Select()
